$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "staff" to "branch"
$ws.Name = "branch"

# Update the defined names (Print_Area / Print_Titles) to reference the new sheet name
$printArea = $wb.Names.Item("branch!Print_Area")
$printArea.RefersTo = "=branch!`$A`$1:`$A`$88"
$printTitles = $wb.Names.Item("branch!Print_Titles")
$printTitles.RefersTo = "=branch!`$1:`$1"

# Populate header row + data rows for columns B, C, D first (matches the
# author's original fill order), then column A (id column) afterwards.
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "location"
$ws.Range("D1").Value = "staffQuota"

$ws.Range("B2").Value = "NTU"
$ws.Range("C2").Value = "North spine Plaza"
$ws.Range("D2").Value = 8

$ws.Range("B3").Value = "JP"
$ws.Range("C3").Value = "Jurong point"
$ws.Range("D3").Value = 15

$ws.Range("B4").Value = "JE"
$ws.Range("C4").Value = "Jurong east"
$ws.Range("D4").Value = 11

$ws.Range("A1").Value = "id"
$ws.Range("A2").Value = "477c0c7e-9d46-4202-969d-f3dd1933a575"
$ws.Range("A3").Value = "67136f7c-fcd0-45f1-8859-9e3d183faeb3"
$ws.Range("A4").Value = "5cc0e578-41b6-4e7d-b6e8-5f287be3e857"

# A2:A4 lose their original style (A1 keeps it)
$ws.Range("A2").Style = "Normal"
$ws.Range("A3").Style = "Normal"
$ws.Range("A4").Style = "Normal"

# Update the frozen-pane selection to match the new data range
$ws.Range("A2:A4").Select()
